# Generate Report for Handoff
#
# The previous handoff round (86805cd3.../f200f13a...) is replaced by a new
# handoff round (e11abc11.../ffff8d28186a...): new source-file tokens, a new
# "Ready for handoff" status, new handoff timestamps, and (since a fresh
# handoff has no handback yet) the Latest Target File / Latest Handback File
# / Latest Handback DateTime columns are cleared out on the language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New values for this handoff round
# ---------------------------------------------------------------------
$oldFile1 = "86805cd3-1a9f-4ce1-8224-5125e6abfa5b.md"
$newFile1 = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"

$oldFile2 = "f200f13a-da38-4319-9c3c-540a262a9d06.md"
$newFile2 = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"

$newStatus       = "Ready for handoff"
$newHandoffDate  = "2016-03-21 12:56:36"

$newZhXlf        = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$newDeXlf        = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"

$newHandoffDatetimeZh = "2016-03-21 12:56:32"
$newHandoffDatetimeDe = "2016-03-21 12:56:36"

$newHandbackDatetime  = "0001-01-01 00:00:00"

# =======================================================================
# Sheet "Overview"
# =======================================================================
$wsOverview = $wb.Worksheets("Overview")

# Plain (non-hyperlink) cells first.
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = $newHandoffDate
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $newHandoffDate

# Hyperlinked file-name cells: drop the old links and add fresh ones so the
# link target/display and the cell text all move together.
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fd3ab389ae77a94f11c4cff8771defcc390b53b8/e2e/$newFile1", "", "", $newFile1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fd3ab389ae77a94f11c4cff8771defcc390b53b8/e2e/$newFile2", "", "", $newFile2)

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$wsZh = $wb.Worksheets("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("E2").Value = $newHandoffDatetimeZh
$wsZh.Range("H2").Value = $newHandbackDatetime

$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("E3").Value = $newHandoffDatetimeZh
$wsZh.Range("H3").Value = $newHandbackDatetime

# No handback yet this round -- clear the "Latest Target File" /
# "Latest Handback File" cells entirely (they disappear, not just go blank).
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7ed7748e6577e86694364847189c0e6725dee5a0/e2e/$newFile1", "", "", $newFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3fe474cbfd3c4f986f871c9dc04951b8f96021ad/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf", "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7ed7748e6577e86694364847189c0e6725dee5a0/e2e/$newFile2", "", "", $newFile2)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3fe474cbfd3c4f986f871c9dc04951b8f96021ad/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf", "", "", $newZhXlf)

# =======================================================================
# Sheet "de-de"
# =======================================================================
$wsDe = $wb.Worksheets("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("E2").Value = $newHandoffDatetimeDe
$wsDe.Range("H2").Value = $newHandbackDatetime

$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("E3").Value = $newHandoffDatetimeDe
$wsDe.Range("H3").Value = $newHandbackDatetime

$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f1aa87184eebf6a7c6d31b8f1a36c5708243271e/e2e/$newFile1", "", "", $newFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3cc86ddeacb2ba8d1d32d954d0abfcdef9452d4b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf", "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f1aa87184eebf6a7c6d31b8f1a36c5708243271e/e2e/$newFile2", "", "", $newFile2)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3cc86ddeacb2ba8d1d32d954d0abfcdef9452d4b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf", "", "", $newDeXlf)
